{"js": "// Replace the multiplication-problem text in each table cell with its\n// updated value. Every \"NNN\u00d7N=\" string in the document is unique, so a\n// plain text search-and-replace for each pair is safe and unambiguous.\nconst replacements = [\n  [\"829\u00d73=\", \"316\u00d73=\"],\n  [\"270\u00d79=\", \"498\u00d78=\"],\n  [\"700\u00d77=\", \"119\u00d78=\"],\n  [\"196\u00d74=\", \"268\u00d77=\"],\n  [\"264\u00d72=\", \"751\u00d75=\"],\n  [\"993\u00d73=\", \"501\u00d76=\"],\n  [\"130\u00d76=\", \"192\u00d78=\"],\n  [\"326\u00d78=\", \"418\u00d78=\"],\n  [\"251\u00d77=\", \"370\u00d75=\"],\n  [\"273\u00d74=\", \"796\u00d77=\"],\n  [\"900\u00d77=\", \"746\u00d79=\"],\n  [\"284\u00d77=\", \"348\u00d75=\"],\n  [\"677\u00d74=\", \"975\u00d74=\"],\n  [\"748\u00d73=\", \"638\u00d79=\"],\n  [\"850\u00d77=\", \"694\u00d72=\"],\n  [\"424\u00d72=\", \"462\u00d73=\"],\n  [\"782\u00d75=\", \"923\u00d75=\"],\n  [\"155\u00d73=\", \"186\u00d72=\"],\n  [\"120\u00d77=\", \"988\u00d72=\"],\n  [\"396\u00d73=\", \"482\u00d73=\"],\n  [\"911\u00d77=\", \"782\u00d73=\"],\n  [\"140\u00d72=\", \"349\u00d79=\"],\n  [\"861\u00d75=\", \"502\u00d78=\"],\n  [\"268\u00d78=\", \"439\u00d79=\"],\n  [\"325\u00d76=\", \"894\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell with its\n# updated value. Every \"NNN\u00d7N=\" string in the document is unique, so a\n# plain Find/Replace for each pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"829\u00d73=\"; new = \"316\u00d73=\"},\n    @{old = \"270\u00d79=\"; new = \"498\u00d78=\"},\n    @{old = \"700\u00d77=\"; new = \"119\u00d78=\"},\n    @{old = \"196\u00d74=\"; new = \"268\u00d77=\"},\n    @{old = \"264\u00d72=\"; new = \"751\u00d75=\"},\n    @{old = \"993\u00d73=\"; new = \"501\u00d76=\"},\n    @{old = \"130\u00d76=\"; new = \"192\u00d78=\"},\n    @{old = \"326\u00d78=\"; new = \"418\u00d78=\"},\n    @{old = \"251\u00d77=\"; new = \"370\u00d75=\"},\n    @{old = \"273\u00d74=\"; new = \"796\u00d77=\"},\n    @{old = \"900\u00d77=\"; new = \"746\u00d79=\"},\n    @{old = \"284\u00d77=\"; new = \"348\u00d75=\"},\n    @{old = \"677\u00d74=\"; new = \"975\u00d74=\"},\n    @{old = \"748\u00d73=\"; new = \"638\u00d79=\"},\n    @{old = \"850\u00d77=\"; new = \"694\u00d72=\"},\n    @{old = \"424\u00d72=\"; new = \"462\u00d73=\"},\n    @{old = \"782\u00d75=\"; new = \"923\u00d75=\"},\n    @{old = \"155\u00d73=\"; new = \"186\u00d72=\"},\n    @{old = \"120\u00d77=\"; new = \"988\u00d72=\"},\n    @{old = \"396\u00d73=\"; new = \"482\u00d73=\"},\n    @{old = \"911\u00d77=\"; new = \"782\u00d73=\"},\n    @{old = \"140\u00d72=\"; new = \"349\u00d79=\"},\n    @{old = \"861\u00d75=\"; new = \"502\u00d78=\"},\n    @{old = \"268\u00d78=\"; new = \"439\u00d79=\"},\n    @{old = \"325\u00d76=\"; new = \"894\u00d76=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null\n}\n"}
